# Apply the "bandit problems" topics edit to the Topics worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows right after the current row 30 ("Bayesian model
# selection", week 12) so the new week-12 topics land right before the
# week-13 "dimensionality reduction" block which is currently rows 31-33.
$insertRange = $ws.Range("A31:A33").EntireRow
$insertRange.Insert()

# New week-12 topics (rows 31-33), formatted like the existing rows
# (B column carries the same explicit black-font style used throughout
# the sheet, matching style index 1 / font rgb FF000000).
$ws.Range("A31").Value = "bandit problems"
$ws.Range("B31").Value = "lecture"
$ws.Range("B31").Font.Color = 0
$ws.Range("C31").Value = 12

$ws.Range("A32").Value = "Thompson sampling"
$ws.Range("B32").Value = "lecture"
$ws.Range("B32").Font.Color = 0
$ws.Range("C32").Value = 12

$ws.Range("A33").Value = "UCB strategy"
$ws.Range("B33").Value = "lecture"
$ws.Range("B33").Font.Color = 0
$ws.Range("C33").Value = 12

# Fix the typo in the final row's topic ("catchup" -> "catch up").
$ws.Range("A37").Value = "catch up and/or advanced topics"

$ws.Range("A37").Select()
